# Auto-generated edit script: update referee stats values and refresh timestamps
$wb = $excel.ActiveWorkbook
$newTimestamp = "2025-10-30 03:03:17"

# ---- Sheet: Главные ----
$ws = $wb.Worksheets.Item("Главные")

$ws.Range("C2").Value = 19
$ws.Range("D2").Value = 466
$ws.Range("E2").Value = 198
$ws.Range("F2").Value = 268
$ws.Range("G2").Value = 24.53
$ws.Range("H2").Value = 10.42
$ws.Range("I2").Value = 14.11
$ws.Range("J2").Value = 84
$ws.Range("K2").Value = 104

$ws.Range("C5").Value = 19
$ws.Range("D5").Value = 315
$ws.Range("E5").Value = 167
$ws.Range("F5").Value = 148
$ws.Range("G5").Value = 16.58
$ws.Range("H5").Value = 8.789999999999999
$ws.Range("I5").Value = 7.79
$ws.Range("J5").Value = 81
$ws.Range("K5").Value = 74

$ws.Range("C15").Value = 13
$ws.Range("D15").Value = 190
$ws.Range("E15").Value = 84
$ws.Range("F15").Value = 106
$ws.Range("G15").Value = 14.62
$ws.Range("H15").Value = 6.46
$ws.Range("I15").Value = 8.15
$ws.Range("J15").Value = 42
$ws.Range("K15").Value = 53

$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 206
$ws.Range("E17").Value = 68
$ws.Range("F17").Value = 138
$ws.Range("G17").Value = 15.85
$ws.Range("H17").Value = 5.23
$ws.Range("I17").Value = 10.62
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = 54

$ws.Range("C20").Value = 18
$ws.Range("D20").Value = 305
$ws.Range("E20").Value = 118
$ws.Range("F20").Value = 187
$ws.Range("G20").Value = 16.94
$ws.Range("H20").Value = 6.56
$ws.Range("I20").Value = 10.39
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = 66
$ws.Range("V20").Value = 12

$ws.Range("C22").Value = 12
$ws.Range("D22").Value = 252
$ws.Range("E22").Value = 94
$ws.Range("F22").Value = 158
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = 7.83
$ws.Range("I22").Value = 13.17
$ws.Range("J22").Value = 47
$ws.Range("K22").Value = 49

$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 308
$ws.Range("E24").Value = 149
$ws.Range("F24").Value = 159
$ws.Range("G24").Value = 16.21
$ws.Range("H24").Value = 7.84
$ws.Range("I24").Value = 8.369999999999999
$ws.Range("J24").Value = 72
$ws.Range("K24").Value = 77
$ws.Range("V24").Value = 10

for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA" + $r).Value = $newTimestamp
}

# ---- Sheet: Линейные ----
$ws = $wb.Worksheets.Item("Линейные")

$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 160
$ws.Range("E4").Value = 62
$ws.Range("F4").Value = 98
$ws.Range("G4").Value = 17.78
$ws.Range("H4").Value = 6.89
$ws.Range("I4").Value = 10.89
$ws.Range("J4").Value = 31
$ws.Range("K4").Value = 39

$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 223
$ws.Range("E6").Value = 92
$ws.Range("F6").Value = 131
$ws.Range("G6").Value = 18.58
$ws.Range("H6").Value = 7.67
$ws.Range("I6").Value = 10.92
$ws.Range("J6").Value = 46
$ws.Range("K6").Value = 63

$ws.Range("C8").Value = 17
$ws.Range("D8").Value = 256
$ws.Range("E8").Value = 104
$ws.Range("F8").Value = 152
$ws.Range("G8").Value = 15.06
$ws.Range("H8").Value = 6.12
$ws.Range("I8").Value = 8.94
$ws.Range("J8").Value = 47
$ws.Range("K8").Value = 61

$ws.Range("C9").Value = 18
$ws.Range("D9").Value = 318
$ws.Range("E9").Value = 146
$ws.Range("F9").Value = 172
$ws.Range("G9").Value = 17.67
$ws.Range("H9").Value = 8.109999999999999
$ws.Range("I9").Value = 9.56
$ws.Range("J9").Value = 68
$ws.Range("K9").Value = 81

$ws.Range("C13").Value = 18
$ws.Range("D13").Value = 306
$ws.Range("E13").Value = 156
$ws.Range("F13").Value = 150
$ws.Range("G13").Value = 17
$ws.Range("H13").Value = 8.67
$ws.Range("I13").Value = 8.33
$ws.Range("J13").Value = 78
$ws.Range("K13").Value = 70
$ws.Range("V13").Value = 10

$ws.Range("C22").Value = 12
$ws.Range("D22").Value = 185
$ws.Range("E22").Value = 84
$ws.Range("F22").Value = 101
$ws.Range("G22").Value = 15.42
$ws.Range("H22").Value = 7
$ws.Range("I22").Value = 8.42
$ws.Range("J22").Value = 42
$ws.Range("K22").Value = 48

for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA" + $r).Value = $newTimestamp
}
